# Generate Report for Handoff
# Update status text + timestamps to reflect handoff, and widen the
# "status"/datetime columns on all three sheets to fit the new text.

$wb = $excel.ActiveWorkbook

$ovw = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Status text: "In Translation" -> "Ready for handoff"
$ovw.Range("E2").Value = "Ready for handoff"
$ovw.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Timestamps bumped forward (zh-cn handoff generated first, then de-de)
$zhcn.Range("H2").Value = "2016-09-02 05:03:17"
$dede.Range("H2").Value = "2016-09-02 05:03:21"
$ovw.Range("G2").Value = "2016-09-02 05:03:21"

# Widen the Status / datetime columns so the longer text fits.
# (Target stored column width is 17.2159881591797 "characters"; the
# ColumnWidth setter here rounds to the nearest whole pixel, so feed it
# the un-rounded character width and let it land on the nearest pixel.)
$targetColumnWidth = 17.2159881591797 - (5 / 6)
$ovw.Columns.Item(5).ColumnWidth = $targetColumnWidth
$ovw.Columns.Item(6).ColumnWidth = $targetColumnWidth
$zhcn.Columns.Item(3).ColumnWidth = $targetColumnWidth
$dede.Columns.Item(3).ColumnWidth = $targetColumnWidth
